$d = $word.ActiveDocument

function ReplaceText($find, $replace) {
    $r = $d.Content
    $found = $r.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $find"
    }
}

# 1) Split the "Método" sentence: insert a line break before "Viaje didática..."
$find1 = "exerc" + [char]0xED + "cios. Viaje"
$repl1 = "exerc" + [char]0xED + "cios. ^lViaje"
ReplaceText $find1 $repl1

# 2) Split the Bibliografia paragraph into one line per numbered reference,
#    by inserting a manual line break before each item number (2 through 10).
ReplaceText "2.BITTON" "^l2.BITTON"
ReplaceText "3.CHERNICHARO" "^l3.CHERNICHARO"
ReplaceText "4.DEZOTTI" "^l4.DEZOTTI"

$find5 = "5.JORD" + [char]0xC3 + "O"
$repl5 = "^l5.JORD" + [char]0xC3 + "O"
ReplaceText $find5 $repl5

ReplaceText "6.METCALF" "^l6.METCALF"
ReplaceText "7.VON SPERLING" "^l7.VON SPERLING"
ReplaceText "8.VON SPERLING" "^l8.VON SPERLING"
ReplaceText "9.VON SPERLING" "^l9.VON SPERLING"
ReplaceText "10.Fugita" "^l10.Fugita"
